$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.144.97"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "3.127.21"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.92%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.123.95"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.475"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("E13").Value = "  -2.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D16").Value = "3.647.46"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "67.056.67"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("D20").Value = "3.129.01"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "489.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.49%  "
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.05%  "
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.54%  "
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "48.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.62%  "
$ws.Range("E38").Value = "  -2.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.312"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.58%  "
$ws.Range("E44").Value = "  +3.37%  "
$ws.Range("D45").Value = "2.802.37"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "376.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.38%  "
$ws.Range("E47").Value = "  -1.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("E51").Value = "  +2.29%  "
